$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -7.770999999999999
    21 = -8.1
    23 = -7.873999999999998
    25 = -8.228999999999999
    53 = -7.507000000000001
    57 = -7.904000000000001
    59 = -8.061
    69 = -7.221000000000001
    79 = -7.772
    83 = -8.199
    93 = -7.517
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
